$wb = $excel.ActiveWorkbook

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/2ed1f06206c1ae22237b60455cf21e13971dfc47/e2e/fa9defb8-9419-411a-8156-ead20580a9f8.md"
$mdName = "fa9defb8-9419-411a-8156-ead20580a9f8.md"
$errMsg = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/3482209a45cf52ebd53ce2b2e3a844ff78d80f38/e2e/fa9defb8-9419-411a-8156-ead20580a9f8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/2ed1f06206c1ae22237b60455cf21e13971dfc47/e2e/fa9defb8-9419-411a-8156-ead20580a9f8.md.'

# ---------- zh-cn sheet ----------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen "Latest Target File", "Latest Handback File" and "Error Detail" columns to 40
$wsZh.Columns.Item(10).ColumnWidth = 40
$wsZh.Columns.Item(11).ColumnWidth = 40
$wsZh.Columns.Item(18).ColumnWidth = 40

# J2: Latest Target File -> link to the handback md file
$wsZh.Range("J2").Value = $mdName
$wsZh.Hyperlinks.Add($wsZh.Range("J2"), $mdUrl, "", "", $mdName)

# K2: Latest Handback File -> zh-cn xlf file name
$wsZh.Range("K2").Value = "fa9defb8-9419-411a-8156-ead20580a9f8.3198b0bd91bf75f8d55dfd9a7097f656dbdbd338.zh-cn.xlf"

# L2: Latest Handback DateTime -> refreshed timestamp
$wsZh.Range("L2").Value = "2018-03-08 19:44:23"

# R2: Error Detail -> version mismatch message
$wsZh.Range("R2").Value = $errMsg

# ---------- de-de sheet ----------
$wsDe = $wb.Worksheets.Item("de-de")

# Widen "Latest Target File", "Latest Handback File" and "Error Detail" columns to 40
$wsDe.Columns.Item(10).ColumnWidth = 40
$wsDe.Columns.Item(11).ColumnWidth = 40
$wsDe.Columns.Item(18).ColumnWidth = 40

# J2: Latest Target File -> link to the handback md file
$wsDe.Range("J2").Value = $mdName
$wsDe.Hyperlinks.Add($wsDe.Range("J2"), $mdUrl, "", "", $mdName)

# K2: Latest Handback File -> de-de xlf file name
$wsDe.Range("K2").Value = "fa9defb8-9419-411a-8156-ead20580a9f8.3198b0bd91bf75f8d55dfd9a7097f656dbdbd338.de-de.xlf"

# L2: Latest Handback DateTime -> new timestamp (distinct from zh-cn's)
$wsDe.Range("L2").Value = "2018-03-08 19:44:45"

# R2: Error Detail -> version mismatch message (same text as zh-cn)
$wsDe.Range("R2").Value = $errMsg
